# Added choose tract weighted by age group
# - Swap the GEO.id / GEO.id2 / display-label columns so that the
#   human-readable label is first, the GEO.id text second and the GEO.id2
#   numeric code third.
# - Rename the age-bucket headers (drop the leading "T " prefix).
# - Append a "Total" row (58) that sums each age-bucket / total column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Rotate columns A/B/C (GEO.id, GEO.id2, display-label) for every
#    data row. Before: A=GEO.id (text), B=GEO.id2 (number), C=label
#    (text). After: A=label (text), B=GEO.id (text), C=GEO.id2 (number).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 57; $r++) {
    $oldA = $ws.Cells.Item($r, 1).Text
    $oldB = $ws.Cells.Item($r, 2).Value2
    $oldC = $ws.Cells.Item($r, 3).Text

    $ws.Cells.Item($r, 1).Value = $oldC
    $ws.Cells.Item($r, 2).Value = $oldA
    $ws.Cells.Item($r, 3).Value = $oldB
}

# ---------------------------------------------------------------------
# 2. Add the "Total" row at the bottom (row 58) with SUM formulas over
#    the 56 data rows (rows 2-57).
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Total"
$ws.Range("B58").Value = "Total"
$ws.Range("C58").Value = "Total"

$sumCols = @("D", "E", "F", "G", "H", "I", "J", "K", "L")
foreach ($col in $sumCols) {
    $ws.Range($col + "58").Formula = "=SUM(" + $col + "2:" + $col + "57)"
}

# ---------------------------------------------------------------------
# 3. Re-label the header row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "GEO.display-label"
$ws.Range("B1").Value = "GEO.id"
$ws.Range("C1").Value = "GEO.id2"
$ws.Range("D1").Value = "0-14"
$ws.Range("E1").Value = "15-19"
$ws.Range("F1").Value = "20-24"
$ws.Range("G1").Value = "25-34"
$ws.Range("H1").Value = "35-44"
$ws.Range("I1").Value = "45-54"
$ws.Range("J1").Value = "55-64"
$ws.Range("K1").Value = "65+"
$ws.Range("L1").Value = "ALL TOTAL"

# ---------------------------------------------------------------------
# 4. Cosmetic bits: widen column D a touch (it now holds the age-bucket
#    headers) and move the active selection like the authored workbook.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 8.6
$ws.Range("M4").Select()
